# Adds three new character styles (GaNStyle, GaNParagraph, GaNLinks) to the
# document's style sheet and applies them to the relevant runs, matching the
# commit "Add styles to the new paragraphs".

$d = $word.ActiveDocument

# --- 1. Create the character styles -----------------------------------

$gaNStyle = $d.Styles.Add("GaNStyle", 2)            # wdStyleTypeCharacter
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)    # wdStyleTypeCharacter
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)            # wdStyleTypeCharacter
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608                      # 0x800000 -> OOXML 000080 (navy)
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1                        # wdUnderlineSingle

# --- 2. Apply GaNStyle to every "Dates à utiliser..." run -------------

$datesText = "Dates à utiliser pour la Campagne 2022 Constellation du Taureau: 16-25 janvier"
$rng = $d.Content
while ($rng.Find.Execute($datesText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNStyle"
    $rng.Collapse(0)
}

# --- 3. Apply GaNParagraph to the campaign description paragraph ------

$quote = [char]0x2019
$paragraphText = "Vous allez participer à une campagne mondiale d" + $quote + "observation pour détecter les plus faibles étoiles visibles afin de mesurer la pollution lumineuse sur un site donné. Partout dans le monde, en localisant et en observant la Constellation du Taureau dans le ciel nocturne et en la comparant aux cartes stellaires, les participants, apprendront comment l" + $quote + "éclairage, dans leur environnement local, influence la pollution lumineuse. Vos contributions à la base de données en ligne permettront de mesurer la qualité du ciel nocturne."

$rng = $d.Content
if ($rng.Find.Execute($paragraphText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNParagraph"
}

# --- 4. Apply GaNLinks to the credits/links run ------------------------

$linksText = "Les cartes figurant dans ce document ont été établies par Jenik Hollan, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$rng = $d.Content
if ($rng.Find.Execute($linksText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $rng.Style = "GaNLinks"
}
